$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1199.2084
$ws.Range("J19").Value = 1163
$ws.Range("L19").Value = 1163
$ws.Range("N19").Value = -1513
$ws.Range("H40").Value = 12352057
$ws.Range("I40").Value = 4599.75
$ws.Range("J40").Value = 22230022
$ws.Range("K40").Value = 4599.75
$ws.Range("L40").Value = 22230022
$ws.Range("M40").Value = -4424.75
$ws.Range("N40").Value = -22230372
$ws.Range("H42").Value = 5649
$ws.Range("H58").Value = 3169.4
$ws.Range("J58").Value = 10004.667
$ws.Range("L58").Value = 30014.001
$ws.Range("N58").Value = -30314.001
$ws.Range("H62").Value = 2712.6667
$ws.Range("I62").Value = 2129.8572
$ws.Range("K62").Value = 2129.8572
$ws.Range("M62").Value = -1505.8572
$ws.Range("H65").Value = 2712.6667
$ws.Range("I65").Value = 2129.8572
$ws.Range("K65").Value = 10649.286
$ws.Range("M65").Value = -7529.286
$ws.Range("H70").Value = 7187.3335
$ws.Range("I70").Value = 1424.7
$ws.Range("K70").Value = 4274.1
$ws.Range("M70").Value = -4004.1
$ws.Range("H73").Value = 7187.3335
$ws.Range("I73").Value = 1424.7
$ws.Range("K73").Value = 4274.1
$ws.Range("M73").Value = -3338.1
$ws.Range("H86").Value = 13685.363
$ws.Range("I86").Value = 16385.428
$ws.Range("J86").Value = 8960.25
$ws.Range("K86").Value = 16385.428
$ws.Range("L86").Value = 8960.25
$ws.Range("M86").Value = -15262.428
$ws.Range("N86").Value = -11206.25
$ws.Range("H89").Value = 13685.363
$ws.Range("I89").Value = 16385.428
$ws.Range("J89").Value = 8960.25
$ws.Range("K89").Value = 81927.14
$ws.Range("L89").Value = 44801.25
$ws.Range("M89").Value = -76311.14
$ws.Range("N89").Value = -56033.25
$ws.Range("H97").Value = 3209.25
$ws.Range("I97").Value = 294.5
$ws.Range("J97").Value = 6124
$ws.Range("K97").Value = 883.5
$ws.Range("L97").Value = 18372
$ws.Range("M97").Value = -387.5
$ws.Range("N97").Value = -19364
$ws.Range("H98").Value = 1522.6364
$ws.Range("I98").Value = 1561.6
$ws.Range("K98").Value = 1561.6
$ws.Range("M98").Value = -63.59999999999991
$ws.Range("H100").Value = 1491.6666
$ws.Range("I100").Value = 737.5
$ws.Range("K100").Value = 737.5
$ws.Range("M100").Value = -196.5
$ws.Range("H103").Value = 1162.25
$ws.Range("J103").Value = 1749.5
$ws.Range("L103").Value = 5248.5
$ws.Range("N103").Value = -6420.5
$ws.Range("H111").Value = 2129
$ws.Range("I111").Value = 2129
$ws.Range("K111").Value = 6387
$ws.Range("M111").Value = -3320
$ws.Range("H122").Value = 1522.6364
$ws.Range("I122").Value = 1561.6
$ws.Range("K122").Value = 4684.799999999999
$ws.Range("M122").Value = -2234.799999999999
$ws.Range("H123").Value = 300000.5
$ws.Range("J123").Value = 300000.5
$ws.Range("L123").Value = 300000.5
$ws.Range("N123").Value = -309800.5
$ws.Range("H132").Value = 2934.111
$ws.Range("I132").Value = 2675.875
$ws.Range("K132").Value = 8027.625
$ws.Range("M132").Value = -5497.625
$ws.Range("H137").Value = 3219.76
$ws.Range("I137").Value = 2605.182
$ws.Range("J137").Value = 3702.6428
$ws.Range("K137").Value = 7815.545999999999
$ws.Range("L137").Value = 11107.9284
$ws.Range("M137").Value = -5265.545999999999
$ws.Range("N137").Value = -16207.9284
$ws.Range("H138").Value = 4927.5576
$ws.Range("I138").Value = 3765.8235
$ws.Range("J138").Value = 5376.409
$ws.Range("K138").Value = 11297.4705
$ws.Range("L138").Value = 16129.227
$ws.Range("M138").Value = -6157.470499999999
$ws.Range("N138").Value = -26409.227
$ws.Range("H141").Value = 4515.9
$ws.Range("I141").Value = 4672.6665
$ws.Range("K141").Value = 14017.9995
$ws.Range("M141").Value = -8837.999500000002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2392.0784
$ws.Range("I32").Value = 2234.3264
$ws.Range("J32").Value = 6257
$ws.Range("K32").Value = 2234.3264
$ws.Range("L32").Value = 6257
$ws.Range("M32").Value = -1947.3264
$ws.Range("N32").Value = -6831
$ws.Range("H61").Value = 200002160
$ws.Range("I61").Value = 200002160
$ws.Range("K61").Value = 200002160
$ws.Range("M61").Value = -200001948
$ws.Range("H74").Value = 20411268
$ws.Range("I74").Value = 30306028
$ws.Range("J74").Value = 3328.5
$ws.Range("K74").Value = 30306028
$ws.Range("L74").Value = 3328.5
$ws.Range("M74").Value = -30305154
$ws.Range("N74").Value = -5076.5
$ws.Range("H77").Value = 20411268
$ws.Range("I77").Value = 30306028
$ws.Range("J77").Value = 3328.5
$ws.Range("K77").Value = 151530140
$ws.Range("L77").Value = 16642.5
$ws.Range("M77").Value = -151525772
$ws.Range("N77").Value = -25378.5
$ws.Range("H102").Value = 5264635
$ws.Range("I102").Value = 5883863
$ws.Range("J102").Value = 1199.5
$ws.Range("K102").Value = 5883863
$ws.Range("L102").Value = 1199.5
$ws.Range("M102").Value = -5882241
$ws.Range("N102").Value = -4443.5
$ws.Range("H122").Value = 6400.4814
$ws.Range("I122").Value = 5230.2
$ws.Range("J122").Value = 9744.143
$ws.Range("K122").Value = 15690.6
$ws.Range("L122").Value = 29232.429
$ws.Range("M122").Value = -13240.6
$ws.Range("N122").Value = -34132.429
$ws.Range("H132").Value = 4874538
$ws.Range("I132").Value = 2705486
$ws.Range("K132").Value = 8116458
$ws.Range("M132").Value = -8113928
$ws.Range("H136").Value = 200002160
$ws.Range("I136").Value = 200002160
$ws.Range("K136").Value = 600006480
$ws.Range("M136").Value = -600003930

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 599.5
$ws.Range("I10").Value = 599.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 599.5
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -459.5
$ws.Range("H86").Value = 3983
$ws.Range("I86").Value = 3976.2
$ws.Range("K86").Value = 3976.2
$ws.Range("M86").Value = -2853.2
$ws.Range("H89").Value = 3983
$ws.Range("I89").Value = 3976.2
$ws.Range("K89").Value = 19881
$ws.Range("M89").Value = -14265
$ws.Range("H99").Value = 1603.2693
$ws.Range("I99").Value = 1566.1428
$ws.Range("J99").Value = 1759.2
$ws.Range("K99").Value = 1566.1428
$ws.Range("L99").Value = 1759.2
$ws.Range("M99").Value = -68.14280000000008
$ws.Range("N99").Value = -4755.2
$ws.Range("H107").Value = 49206.047
$ws.Range("I107").Value = 1592.9412
$ws.Range("K107").Value = 1592.9412
$ws.Range("M107").Value = 327.0588
$ws.Range("H134").Value = 13558440
$ws.Range("I134").Value = 13924855
$ws.Range("J134").Value = 1098
$ws.Range("K134").Value = 41774565
$ws.Range("L134").Value = 3294
$ws.Range("M134").Value = -41772030
$ws.Range("N134").Value = -8364

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5863.421
$ws.Range("I22").Value = 6587.8125
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 6587.8125
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -6237.8125
$ws.Range("N22").Value = -2700
$ws.Range("H31").Value = 2592.6191
$ws.Range("I31").Value = 2517.8572
$ws.Range("J31").Value = 2742.1428
$ws.Range("K31").Value = 2517.8572
$ws.Range("L31").Value = 2742.1428
$ws.Range("M31").Value = -2222.8572
$ws.Range("N31").Value = -3332.1428
$ws.Range("H34").Value = 2592.6191
$ws.Range("I34").Value = 2517.8572
$ws.Range("J34").Value = 2742.1428
$ws.Range("K34").Value = 2517.8572
$ws.Range("L34").Value = 2742.1428
$ws.Range("M34").Value = -2315.8572
$ws.Range("N34").Value = -3146.1428
$ws.Range("H48").Value = 40001
$ws.Range("J48").Value = 40001
$ws.Range("L48").Value = 40001
$ws.Range("N48").Value = -40953
$ws.Range("H58").Value = 16675758
$ws.Range("I58").Value = 22739146
$ws.Range("J58").Value = 1440.625
$ws.Range("K58").Value = 22739146
$ws.Range("L58").Value = 1440.625
$ws.Range("M58").Value = -22738943
$ws.Range("N58").Value = -1846.625
$ws.Range("H59").Value = 75000
$ws.Range("I59").Value = 10000
$ws.Range("J59").Value = 140000
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 140000
$ws.Range("M59").Value = -8855
$ws.Range("N59").Value = -142290
$ws.Range("H86").Value = 9582.833000000001
$ws.Range("I86").Value = 6666
$ws.Range("J86").Value = 12499.667
$ws.Range("K86").Value = 6666
$ws.Range("L86").Value = 12499.667
$ws.Range("M86").Value = -5543
$ws.Range("N86").Value = -14745.667
$ws.Range("H88").Value = 27333.334
$ws.Range("J88").Value = 27333.334
$ws.Range("L88").Value = 27333.334
$ws.Range("N88").Value = -28145.334
$ws.Range("H89").Value = 9582.833000000001
$ws.Range("I89").Value = 6666
$ws.Range("J89").Value = 12499.667
$ws.Range("K89").Value = 33330
$ws.Range("L89").Value = 62498.335
$ws.Range("M89").Value = -27714
$ws.Range("N89").Value = -73730.33499999999
$ws.Range("H91").Value = 27333.334
$ws.Range("J91").Value = 27333.334
$ws.Range("L91").Value = 27333.334
$ws.Range("N91").Value = -30141.334
$ws.Range("H94").Value = 1787.2858
$ws.Range("I94").Value = 2124.2
$ws.Range("K94").Value = 2124.2
$ws.Range("M94").Value = -1673.2
$ws.Range("H122").Value = 3290.4092
$ws.Range("I122").Value = 3294.4375
$ws.Range("J122").Value = 3279.6667
$ws.Range("K122").Value = 9883.3125
$ws.Range("L122").Value = 9839.000100000001
$ws.Range("M122").Value = -7433.3125
$ws.Range("N122").Value = -14739.0001
$ws.Range("H132").Value = 43479750
$ws.Range("I132").Value = 47620456
$ws.Range("J132").Value = 2393.5
$ws.Range("K132").Value = 142861368
$ws.Range("L132").Value = 7180.5
$ws.Range("M132").Value = -142858838
$ws.Range("N132").Value = -12240.5
$ws.Range("H134").Value = 12501900
$ws.Range("J134").Value = 2211.25
$ws.Range("L134").Value = 6633.75
$ws.Range("N134").Value = -11703.75
$ws.Range("H136").Value = 16675758
$ws.Range("I136").Value = 22739146
$ws.Range("J136").Value = 1440.625
$ws.Range("K136").Value = 68217438
$ws.Range("L136").Value = 4321.875
$ws.Range("M136").Value = -68214888
$ws.Range("N136").Value = -9421.875

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 317.9
$ws.Range("J12").Value = 427.42856
$ws.Range("L12").Value = 1282.28568
$ws.Range("N12").Value = -1628.28568
$ws.Range("H56").Value = 14094.8955
$ws.Range("I56").Value = 14094.8955
$ws.Range("K56").Value = 14094.8955
$ws.Range("M56").Value = -13564.8955
$ws.Range("H107").Value = 5108.3335
$ws.Range("J107").Value = 5425
$ws.Range("L107").Value = 16275
$ws.Range("N107").Value = -20115
$ws.Range("H113").Value = 999999
$ws.Range("I113").Value = 999999
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2999997
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -2997827
$ws.Range("H132").Value = 1725
$ws.Range("J132").Value = 1833.3334
$ws.Range("L132").Value = 16500.0006
$ws.Range("N132").Value = -21560.0006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 51399.5
$ws.Range("J15").Value = 51399.5
$ws.Range("L15").Value = 51399.5
$ws.Range("N15").Value = -51975.5
$ws.Range("H70").Value = 8624.75
$ws.Range("I70").Value = 7250
$ws.Range("J70").Value = 9999.5
$ws.Range("K70").Value = 7250
$ws.Range("L70").Value = 9999.5
$ws.Range("M70").Value = -6980
$ws.Range("N70").Value = -10539.5
$ws.Range("H73").Value = 8624.75
$ws.Range("I73").Value = 7250
$ws.Range("J73").Value = 9999.5
$ws.Range("K73").Value = 7250
$ws.Range("L73").Value = 9999.5
$ws.Range("M73").Value = -6314
$ws.Range("N73").Value = -11871.5
$ws.Range("H80").Value = 3448.5652
$ws.Range("I80").Value = 3531.3076
$ws.Range("J80").Value = 3341
$ws.Range("K80").Value = 3531.3076
$ws.Range("L80").Value = 3341
$ws.Range("M80").Value = -2533.3076
$ws.Range("N80").Value = -5337
$ws.Range("H81").Value = 51399.5
$ws.Range("J81").Value = 51399.5
$ws.Range("L81").Value = 51399.5
$ws.Range("N81").Value = -53395.5
$ws.Range("H83").Value = 3448.5652
$ws.Range("I83").Value = 3531.3076
$ws.Range("J83").Value = 3341
$ws.Range("K83").Value = 17656.538
$ws.Range("L83").Value = 16705
$ws.Range("M83").Value = -12664.538
$ws.Range("N83").Value = -26689
$ws.Range("H84").Value = 51399.5
$ws.Range("J84").Value = 51399.5
$ws.Range("L84").Value = 154198.5
$ws.Range("N84").Value = -164182.5
$ws.Range("H97").Value = 4998.6665
$ws.Range("I97").Value = 4998.5
$ws.Range("J97").Value = 4999
$ws.Range("K97").Value = 4998.5
$ws.Range("L97").Value = 4999
$ws.Range("M97").Value = -4502.5
$ws.Range("N97").Value = -5991
$ws.Range("H107").Value = 2202.625
$ws.Range("I107").Value = 1088.2858
$ws.Range("K107").Value = 1088.2858
$ws.Range("M107").Value = 831.7141999999999
$ws.Range("H108").Value = 78000
$ws.Range("J108").Value = 78000
$ws.Range("L108").Value = 78000
$ws.Range("N108").Value = -85680
$ws.Range("H113").Value = 40968.867
$ws.Range("I113").Value = 49545.75
$ws.Range("J113").Value = 6661.3335
$ws.Range("K113").Value = 49545.75
$ws.Range("L113").Value = 6661.3335
$ws.Range("M113").Value = -47375.75
$ws.Range("N113").Value = -11001.3335
$ws.Range("H122").Value = 8231.883
$ws.Range("I122").Value = 5254.6665
$ws.Range("K122").Value = 15763.9995
$ws.Range("M122").Value = -13313.9995
$ws.Range("H132").Value = 2509333
$ws.Range("I132").Value = 2849819.8
$ws.Range("K132").Value = 8549459.399999999
$ws.Range("M132").Value = -8546929.399999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2149.8333
$ws.Range("I7").Value = 2179.8
$ws.Range("K7").Value = 2179.8
$ws.Range("M7").Value = -2067.8
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("K14").Value = 500
$ws.Range("M14").Value = -328
$ws.Range("H20").Value = 50000
$ws.Range("I20").Value = 50000
$ws.Range("K20").Value = 50000
$ws.Range("M20").Value = -49774
$ws.Range("H22").Value = 3236.182
$ws.Range("I22").Value = 3362.3333
$ws.Range("J22").Value = 3084.8
$ws.Range("K22").Value = 3362.3333
$ws.Range("L22").Value = 3084.8
$ws.Range("M22").Value = -3067.3333
$ws.Range("N22").Value = -3674.8
$ws.Range("H27").Value = 3236.182
$ws.Range("I27").Value = 3362.3333
$ws.Range("J27").Value = 3084.8
$ws.Range("K27").Value = 3362.3333
$ws.Range("L27").Value = 3084.8
$ws.Range("M27").Value = -3255.3333
$ws.Range("N27").Value = -3298.8
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = 0
$ws.Range("H53").Value = 25332.666
$ws.Range("I53").Value = 12999.5
$ws.Range("K53").Value = 12999.5
$ws.Range("M53").Value = -12481.5
$ws.Range("H55").Value = 799.2
$ws.Range("I55").Value = 348
$ws.Range("K55").Value = 348
$ws.Range("M55").Value = -175
$ws.Range("H61").Value = 16999.5
$ws.Range("I61").Value = 16999.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 16999.5
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -16797.5
$ws.Range("H113").Value = 16999.5
$ws.Range("I113").Value = 16999.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 16999.5
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -14829.5
$ws.Range("H122").Value = 2996.5386
$ws.Range("I122").Value = 2985.5557
$ws.Range("J122").Value = 3021.25
$ws.Range("K122").Value = 8956.667099999999
$ws.Range("L122").Value = 9063.75
$ws.Range("M122").Value = -6506.667099999999
$ws.Range("N122").Value = -13963.75
$ws.Range("H126").Value = 2149.8333
$ws.Range("I126").Value = 2179.8
$ws.Range("K126").Value = 6539.400000000001
$ws.Range("M126").Value = -4069.400000000001
$ws.Range("H132").Value = 13166583
$ws.Range("I132").Value = 13897874
$ws.Range("K132").Value = 41693622
$ws.Range("M132").Value = -41691092
$ws.Range("H133").Value = 99999
$ws.Range("J133").Value = 99999
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -105059
$ws.Range("H136").Value = 2270.05
$ws.Range("J136").Value = 2500
$ws.Range("L136").Value = 7500
$ws.Range("N136").Value = -12600

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("N22").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("N46").Value = 0
$ws.Range("H62").Value = 7332.5
$ws.Range("J62").Value = 8199
$ws.Range("L62").Value = 8199
$ws.Range("N62").Value = -9447
$ws.Range("H65").Value = 7332.5
$ws.Range("J65").Value = 8199
$ws.Range("L65").Value = 40995
$ws.Range("N65").Value = -47235
$ws.Range("H92").Value = 49833.332
$ws.Range("I92").Value = 50500
$ws.Range("J92").Value = 49500
$ws.Range("K92").Value = 50500
$ws.Range("L92").Value = 49500
$ws.Range("M92").Value = -48004
$ws.Range("N92").Value = -54492
$ws.Range("H107").Value = 1930.1666
$ws.Range("I107").Value = 633.8333
$ws.Range("K107").Value = 1901.4999
$ws.Range("M107").Value = 18.50009999999997
$ws.Range("H110").Value = 220666.33
$ws.Range("J110").Value = 220666.33
$ws.Range("L110").Value = 220666.33
$ws.Range("N110").Value = -228846.33
$ws.Range("H113").Value = 999.6667
$ws.Range("I113").Value = 1299.5
$ws.Range("J113").Value = 400
$ws.Range("K113").Value = 3898.5
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = -1728.5
$ws.Range("N113").Value = -5540
$ws.Range("H126").Value = 1165.3334
$ws.Range("I126").Value = 998.4
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 2995.2
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -525.1999999999998
$ws.Range("N126").Value = -10940
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0

Write-Host "Applied all updates"